$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 161.16  # H15 was 154.91
$ws.Cells.Item(15, 9).Value = 161.16  # I15 was 154.91
$ws.Cells.Item(15, 11).Value = 483.48  # K15 was 464.73
$ws.Cells.Item(15, 13).Value = -314.48  # M15 was -295.73
$ws.Cells.Item(40, 8).Value = 43480390  # H40 was 31251952
$ws.Cells.Item(40, 9).Value = 0  # I40 was 1492.8572
$ws.Cells.Item(40, 10).Value = 43480390  # J40 was 40002080
$ws.Cells.Item(40, 11).Value = 0  # K40 was 1492.8572
$ws.Cells.Item(40, 12).Value = 43480390  # L40 was 40002080
$ws.Cells.Item(40, 13).ClearContents()  # M40 was -1317.8572
$ws.Cells.Item(40, 14).Value = -43480740  # N40 was -40002430
$ws.Cells.Item(116, 8).Value = 4242.095  # H116 was 4162.909
$ws.Cells.Item(116, 10).Value = 3658  # J116 was 3552.7273
$ws.Cells.Item(116, 12).Value = 3658  # L116 was 3552.7273
$ws.Cells.Item(116, 14).Value = -10542  # N116 was -10436.7273
$ws.Cells.Item(129, 8).Value = 713214.7  # H129 was 699791.25
$ws.Cells.Item(129, 9).Value = 341.58334  # I129 was 354.45456
$ws.Cells.Item(129, 10).Value = 927076.6  # J129 was 882977.0600000001
$ws.Cells.Item(129, 11).Value = 1024.75002  # K129 was 1063.36368
$ws.Cells.Item(129, 12).Value = 2781229.8  # L129 was 2648931.18
$ws.Cells.Item(129, 13).Value = 3975.24998  # M129 was 3936.63632
$ws.Cells.Item(129, 14).Value = -2791229.8  # N129 was -2658931.18
$ws.Cells.Item(132, 8).Value = 2030.7142  # H132 was 2081.1765
$ws.Cells.Item(132, 9).Value = 2114.394  # I132 was 2120
$ws.Cells.Item(132, 10).Value = 650  # J132 was 800
$ws.Cells.Item(132, 11).Value = 6343.181999999999  # K132 was 6360
$ws.Cells.Item(132, 12).Value = 1950  # L132 was 2400
$ws.Cells.Item(132, 13).Value = -3813.181999999999  # M132 was -3830
$ws.Cells.Item(132, 14).Value = -7010  # N132 was -7460
$ws.Cells.Item(137, 8).Value = 932.14813  # H137 was 884
$ws.Cells.Item(137, 9).Value = 739.1818  # I137 was 724.43475
$ws.Cells.Item(137, 10).Value = 1781.2  # J137 was 1801.5
$ws.Cells.Item(137, 11).Value = 2217.5454  # K137 was 2173.30425
$ws.Cells.Item(137, 12).Value = 5343.6  # L137 was 5404.5
$ws.Cells.Item(137, 13).Value = 332.4546  # M137 was 376.6957499999999
$ws.Cells.Item(137, 14).Value = -10443.6  # N137 was -10504.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 38186.5  # H2 was 20230.283
$ws.Cells.Item(2, 9).Value = 54313  # I2 was 23521.137
$ws.Cells.Item(2, 11).Value = 54313  # K2 was 23521.137
$ws.Cells.Item(2, 13).Value = -54200  # M2 was -23408.137
$ws.Cells.Item(60, 8).Value = 20000  # H60 was 0
$ws.Cells.Item(60, 9).Value = 20000  # I60 was 0
$ws.Cells.Item(60, 11).Value = 20000  # K60 was 0
$ws.Cells.Item(60, 13).Value = -19267  # M60 was None
$ws.Cells.Item(116, 8).Value = 38186.5  # H116 was 20230.283
$ws.Cells.Item(116, 9).Value = 54313  # I116 was 23521.137
$ws.Cells.Item(116, 11).Value = 54313  # K116 was 23521.137
$ws.Cells.Item(116, 13).Value = -52019  # M116 was -21227.137
$ws.Cells.Item(124, 8).Value = 24720.154  # H124 was 25280.166
$ws.Cells.Item(124, 10).Value = 24720.154  # J124 was 25280.166
$ws.Cells.Item(124, 12).Value = 24720.154  # L124 was 25280.166
$ws.Cells.Item(124, 14).Value = -34540.15399999999  # N124 was -35100.166
$ws.Cells.Item(125, 8).Value = 20349.928  # H125 was 20615.309
$ws.Cells.Item(125, 10).Value = 20349.928  # J125 was 20615.309
$ws.Cells.Item(125, 12).Value = 20349.928  # L125 was 20615.309
$ws.Cells.Item(125, 14).Value = -30189.928  # N125 was -30455.309
$ws.Cells.Item(132, 8).Value = 1325.6471  # H132 was 1457.5238
$ws.Cells.Item(132, 9).Value = 1006.65  # I132 was 1106.7878
$ws.Cells.Item(132, 10).Value = 2485.6365  # J132 was 2743.5557
$ws.Cells.Item(132, 11).Value = 3019.95  # K132 was 3320.3634
$ws.Cells.Item(132, 12).Value = 7456.9095  # L132 was 8230.667099999999
$ws.Cells.Item(132, 13).Value = -489.9499999999998  # M132 was -790.3634000000002
$ws.Cells.Item(132, 14).Value = -12516.9095  # N132 was -13290.6671

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 38186.5  # H3 was 20230.283
$ws.Cells.Item(3, 9).Value = 54313  # I3 was 23521.137
$ws.Cells.Item(3, 11).Value = 54313  # K3 was 23521.137
$ws.Cells.Item(3, 13).Value = -54199  # M3 was -23407.137
$ws.Cells.Item(26, 8).Value = 11367.75  # H26 was 5145.1665
$ws.Cells.Item(26, 9).Value = 3490.3333  # I26 was 5145.1665
$ws.Cells.Item(26, 10).Value = 35000  # J26 was 0
$ws.Cells.Item(26, 11).Value = 3490.3333  # K26 was 5145.1665
$ws.Cells.Item(26, 12).Value = 35000  # L26 was 0
$ws.Cells.Item(26, 13).Value = -3198.3333  # M26 was -4853.1665
$ws.Cells.Item(26, 14).Value = -35584  # N26 was None
$ws.Cells.Item(96, 8).Value = 11233.333  # H96 was 27482.6
$ws.Cells.Item(96, 9).Value = 11233.333  # I96 was 10000
$ws.Cells.Item(96, 10).Value = 0  # J96 was 31853.25
$ws.Cells.Item(96, 11).Value = 11233.333  # K96 was 10000
$ws.Cells.Item(96, 12).Value = 0  # L96 was 31853.25
$ws.Cells.Item(96, 13).Value = -8487.333000000001  # M96 was -7254
$ws.Cells.Item(96, 14).ClearContents()  # N96 was -37345.25
$ws.Cells.Item(107, 8).Value = 7482.65  # H107 was 6289.7085
$ws.Cells.Item(107, 9).Value = 803.5333000000001  # I107 was 702.7895
$ws.Cells.Item(107, 11).Value = 803.5333000000001  # K107 was 702.7895
$ws.Cells.Item(107, 13).Value = 1116.4667  # M107 was 1217.2105
$ws.Cells.Item(134, 8).Value = 3161.5  # H134 was 28991.736
$ws.Cells.Item(134, 9).Value = 2820.111  # I134 was 2304.12
$ws.Cells.Item(134, 10).Value = 3673.5833  # J134 was 80314.08
$ws.Cells.Item(134, 11).Value = 8460.332999999999  # K134 was 6912.36
$ws.Cells.Item(134, 12).Value = 11020.7499  # L134 was 240942.24
$ws.Cells.Item(134, 13).Value = -5925.332999999999  # M134 was -4377.36
$ws.Cells.Item(134, 14).Value = -16090.7499  # N134 was -246012.24

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 35332.29  # H31 was 40431.445
$ws.Cells.Item(31, 9).Value = 3324.7896  # I31 was 3835.5625
$ws.Cells.Item(31, 10).Value = 86010.836  # J31 was 93661.82000000001
$ws.Cells.Item(31, 11).Value = 3324.7896  # K31 was 3835.5625
$ws.Cells.Item(31, 12).Value = 86010.836  # L31 was 93661.82000000001
$ws.Cells.Item(31, 13).Value = -3029.7896  # M31 was -3540.5625
$ws.Cells.Item(31, 14).Value = -86600.836  # N31 was -94251.82000000001
$ws.Cells.Item(34, 8).Value = 35332.29  # H34 was 40431.445
$ws.Cells.Item(34, 9).Value = 3324.7896  # I34 was 3835.5625
$ws.Cells.Item(34, 10).Value = 86010.836  # J34 was 93661.82000000001
$ws.Cells.Item(34, 11).Value = 3324.7896  # K34 was 3835.5625
$ws.Cells.Item(34, 12).Value = 86010.836  # L34 was 93661.82000000001
$ws.Cells.Item(34, 13).Value = -3122.7896  # M34 was -3633.5625
$ws.Cells.Item(34, 14).Value = -86414.836  # N34 was -94065.82000000001
$ws.Cells.Item(63, 8).Value = 24750  # H63 was 0
$ws.Cells.Item(63, 10).Value = 24750  # J63 was 0
$ws.Cells.Item(63, 12).Value = 24750  # L63 was 0
$ws.Cells.Item(63, 14).Value = -26122  # N63 was None
$ws.Cells.Item(66, 8).Value = 24750  # H66 was 0
$ws.Cells.Item(66, 10).Value = 24750  # J66 was 0
$ws.Cells.Item(66, 12).Value = 74250  # L66 was 0
$ws.Cells.Item(66, 14).Value = -81114  # N66 was None
$ws.Cells.Item(94, 8).Value = 7970.1665  # H94 was 8230
$ws.Cells.Item(94, 10).Value = 9124.223  # J94 was 9625.75
$ws.Cells.Item(94, 12).Value = 9124.223  # L94 was 9625.75
$ws.Cells.Item(94, 14).Value = -10026.223  # N94 was -10527.75
$ws.Cells.Item(107, 8).Value = 489.1875  # H107 was 458.22223
$ws.Cells.Item(107, 9).Value = 692.5  # I107 was 521.8333
$ws.Cells.Item(107, 10).Value = 421.41666  # J107 was 426.41666
$ws.Cells.Item(107, 11).Value = 692.5  # K107 was 521.8333
$ws.Cells.Item(107, 12).Value = 421.41666  # L107 was 426.41666
$ws.Cells.Item(107, 13).Value = 1227.5  # M107 was 1398.1667
$ws.Cells.Item(107, 14).Value = -4261.41666  # N107 was -4266.41666
$ws.Cells.Item(134, 8).Value = 17858560  # H134 was 18519958
$ws.Cells.Item(134, 9).Value = 1455.3684  # I134 was 1465.5264
$ws.Cells.Item(134, 10).Value = 55556892  # J134 was 62501380
$ws.Cells.Item(134, 11).Value = 4366.1052  # K134 was 4396.5792
$ws.Cells.Item(134, 12).Value = 166670676  # L134 was 187504140
$ws.Cells.Item(134, 13).Value = -1831.1052  # M134 was -1861.5792
$ws.Cells.Item(134, 14).Value = -166675746  # N134 was -187509210

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(110, 8).Value = 10748.75  # H110 was 11276.637
$ws.Cells.Item(110, 9).Value = 4597  # I110 was 4510.75
$ws.Cells.Item(110, 11).Value = 13791  # K110 was 13532.25
$ws.Cells.Item(110, 13).Value = -9701  # M110 was -9442.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2977.9062  # H132 was 3445.5386
$ws.Cells.Item(132, 9).Value = 2666.7273  # I132 was 3309.9375
$ws.Cells.Item(132, 11).Value = 8000.1819  # K132 was 9929.8125
$ws.Cells.Item(132, 13).Value = -5470.1819  # M132 was -7399.8125
$ws.Cells.Item(134, 8).Value = 13104.23  # H134 was 13311.923
$ws.Cells.Item(134, 10).Value = 13104.23  # J134 was 13311.923
$ws.Cells.Item(134, 12).Value = 39312.69  # L134 was 39935.769
$ws.Cells.Item(134, 14).Value = -44382.69  # N134 was -45005.769

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 1934.3396  # H132 was 1869.9822
$ws.Cells.Item(132, 9).Value = 1812.359  # I132 was 1735.262
$ws.Cells.Item(132, 11).Value = 5437.076999999999  # K132 was 5205.786
$ws.Cells.Item(132, 13).Value = -2907.076999999999  # M132 was -2675.786

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 9933.333000000001  # H54 was 9491.666999999999
$ws.Cells.Item(54, 9).Value = 9900  # I54 was 10000
$ws.Cells.Item(54, 10).Value = 10000  # J54 was 9237.5
$ws.Cells.Item(54, 11).Value = 9900  # K54 was 10000
$ws.Cells.Item(54, 12).Value = 10000  # L54 was 9237.5
$ws.Cells.Item(54, 13).Value = -9380  # M54 was -9480
$ws.Cells.Item(54, 14).Value = -11040  # N54 was -10277.5
$ws.Cells.Item(132, 8).Value = 876.5454999999999  # H132 was 896.86365
$ws.Cells.Item(132, 9).Value = 719.4  # I132 was 786.65625
$ws.Cells.Item(132, 10).Value = 1487.6666  # J132 was 1190.75
$ws.Cells.Item(132, 11).Value = 2158.2  # K132 was 2359.96875
$ws.Cells.Item(132, 12).Value = 4462.9998  # L132 was 3572.25
$ws.Cells.Item(132, 13).Value = 371.8000000000002  # M132 was 170.03125
$ws.Cells.Item(132, 14).Value = -9522.9998  # N132 was -8632.25
